$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 and Row 4 height changes
$ws.Rows.Item(1).RowHeight = 57
$ws.Rows.Item(4).RowHeight = 16.5

# Add the new 2023 column (T) of data, copying the formatting from the
# corresponding 2022 column (S) on each row, then overwriting with the
# new values.
$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("T5").Value = 43.1

$ws.Range("S6").Copy($ws.Range("T6"))
$ws.Range("T6").Value = 19.7

$ws.Range("S7").Copy($ws.Range("T7"))
$ws.Range("T7").Value = 7.8

$ws.Range("S8").Copy($ws.Range("T8"))
$ws.Range("T8").Value = 15.6

# Reset the view selection back to the top-left cell (the authored sheet no
# longer carries a stray selection on Y14).
[void]$ws.Range("A1").Select()
$excel.CutCopyMode = $false
